$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values (columns B:E)
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 values (columns B:E)
$ws.Range("B2").Value = 0.66045457039765831
$ws.Range("C2").Value = 1.6617924319921875
$ws.Range("D2").Value = 0.98377283629348966
$ws.Range("E2").Value = 1.4749013336491272

# Row 3 values (columns B:E)
$ws.Range("B3").Value = 1.8271309109788068
$ws.Range("C3").Value = 0.77217448054328564
$ws.Range("D3").Value = 1.2538360540595634
$ws.Range("E3").Value = 0.80953289967086539

# Update selection to match new range B1:E3
$ws.Range("B1:E3").Select()
